$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")
$ws.Rows(45).Insert()
$ws.Range("R45").Value = "axis"
$ws.Range("S45").Value = "2024-09-20 06:57:43"
